$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.412.42'
$ws.Range('E2').Value = '  +3.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.060.92'
$ws.Range('E3').Value = '  +1.74%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '549.38'
$ws.Range('E5').Value = '  +2.49%  '
$ws.Range('E6').Value = '  +3.49%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.053.50'
$ws.Range('E8').Value = '  +1.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.503'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.47'
$ws.Range('E10').Value = '  +5.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.151'
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.454'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.88'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.556.21'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.279.27'
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.059.35'
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.75'
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '482.35'
$ws.Range('E20').Value = '  +3.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.67'
$ws.Range('E21').Value = '  +2.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.674'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.24'
$ws.Range('E23').Value = '  +3.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.63'
$ws.Range('E24').Value = '  +1.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.63'
$ws.Range('E25').Value = '  +4.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  +2.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.94'
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.99'
$ws.Range('E29').Value = '  +3.95%  '
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.07'
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.46'
$ws.Range('E33').Value = '  +7.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.72'
$ws.Range('E34').Value = '  +3.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '55.50'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.99'
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '469.86'
$ws.Range('E37').Value = '  +2.08%  '
$ws.Range('E38').Value = '  +3.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0398'
$ws.Range('E39').Value = '  +2.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.073.65'
$ws.Range('E40').Value = '  -4.28%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  +0.96%  '
$ws.Range('E43').Value = '  +4.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '28.16'
$ws.Range('E44').Value = '  +1.02%  '
$ws.Range('E45').Value = '  +2.84%  '
$ws.Range('E47').Value = '  +2.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.110'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '116.92'
$ws.Range('E49').Value = '  -2.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0510'
$ws.Range('E50').Value = '  +2.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.07'
$ws.Range('E51').Value = '  +2.53%  '
